# The document's single section has a title-page header/footer pair
# (wdHeaderFooterFirstPage, index 2) and a default header/footer pair
# (wdHeaderFooterPrimary, index 1). Each of the four stories holds one
# inline picture whose display "name" (wp:docPr/@name, mirrored onto
# pic:cNvPr/@name) needs to be swapped:
#   - the two Pearson-logo pictures (footers): image1.png -> image2.png
#   - the two BTEC-logo pictures (headers):    image2.jpg -> image1.jpg
#
# Renaming an InlineShape that isn't the very first thing in its story
# range (true for both footers, which each have several paragraphs
# before the picture) trips a stale-handle error if we set .Name
# straight off a freshly fetched InlineShape. Routing the rename
# through the Selection (select the shape's range, then rename via
# Selection.InlineShapes) sidesteps that and works uniformly for every
# header/footer here.

$d = $word.ActiveDocument
$sec = $d.Sections.First

function Rename-InlineLogo($range, [string]$newName) {
    $range.Select()
    $word.Selection.InlineShapes.Item(1).Name = $newName
}

# Footers (Pearson Edexcel logo): image1.png -> image2.png
Rename-InlineLogo $sec.Footers.Item(2).Range.InlineShapes.Item(1).Range "image2.png"   # footer1.xml (first page), docPr id="2"
Rename-InlineLogo $sec.Footers.Item(1).Range.InlineShapes.Item(1).Range "image2.png"   # footer2.xml (default), docPr id="4"

# Headers (BTEC logo): image2.jpg -> image1.jpg
Rename-InlineLogo $sec.Headers.Item(2).Range.InlineShapes.Item(1).Range "image1.jpg"   # header1.xml (first page), docPr id="1"
Rename-InlineLogo $sec.Headers.Item(1).Range.InlineShapes.Item(1).Range "image1.jpg"   # header2.xml (default), docPr id="3"
